$d = $word.ActiveDocument

# 1. Remove "However, " and capitalize "strict" -> "Strict"
$d.Content.Find.Execute("collected. However, strict principles", $true, $false, $false, $false, $false, $true, 1, $false, "collected. Strict principles", 2)

# 2. Change "24 months" to "5 years" in the retention sentence
$d.Content.Find.Execute("a maximum of 24 months,", $true, $false, $false, $false, $false, $true, 1, $false, "a maximum of 5 years,", 2)
